$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 299 (pushes old rows 299-354 down to 300-355)
$ws.Rows.Item(299).Insert()

# Populate the new row 299. Most fields are copied from the template row (now at 300,
# which holds the data that used to be in row 299), with D, K, L, M, P updated to new values.
$ws.Cells.Item(299, 1).Value = 3
$ws.Cells.Item(299, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(299, 3).Value = "Coquimbo"
$ws.Cells.Item(299, 4).Value = 44694
$ws.Cells.Item(299, 5).Value = 5
$ws.Cells.Item(299, 6).Value = 100112040
$ws.Cells.Item(299, 7).Value = "Cilantro"
$ws.Cells.Item(299, 8).Value = "Sin especificar"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 130
$ws.Cells.Item(299, 11).Value = 3000
$ws.Cells.Item(299, 12).Value = 3300
$ws.Cells.Item(299, 13).Value = 3138
$ws.Cells.Item(299, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(299, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(299, 16).Value = 1046
$ws.Cells.Item(299, 17).Value = 3
$ws.Cells.Item(299, 18).Value = "Hortaliza"

# Copy the date cell's number format/style from the neighboring date cell (row 300)
$ws.Cells.Item(300, 4).Copy()
$ws.Cells.Item(299, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false
